$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not auto-coerced to a number) by using
# Excel's leading-apostrophe quote-prefix, then strip the resulting cell
# style back down to Normal so no extra formatting is left behind.
function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# New board revision: rename the sheet to reflect the fixed PCB revision.
$ws.Name = "BOM_Keys_fixed_PCB2_1_2024-11-0"

# A new BOM line (CHERRY MX1A-11NW switch) was added as the new first
# component row. Insert a row right under the header and shift the
# existing component rows down.
$ws.Rows.Item(2).Insert()

Set-TextValue "A2" "1"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = "MX1A-11NW"
$ws.Range("D2").Value = "BC1,BC2,BC3,BC4"
$ws.Range("E2").Value = "SW-TH_MX1A-11NW"
Set-TextValue "F2" ""
$ws.Range("G2").Value = "MX1A-11NW"
$ws.Range("H2").Value = "CHERRY"
$ws.Range("I2").Value = "C3316924"
$ws.Range("J2").Value = "LCSC"

# The old "SMD_KEYBOARD-SW" line (previously row 5, No.=4) no longer
# belongs on the new board revision - it is now row 6 after the insert
# above, so remove it entirely.
$ws.Rows.Item(6).Delete()

# Renumber the "No." column for the rows that shifted down, so the BOM
# stays sequential (1 = new switch row, 2-4 = the pre-existing parts).
Set-TextValue "A3" "2"
Set-TextValue "A4" "3"
Set-TextValue "A5" "4"

Write-Host "applied board revision update"
